$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-CellText "D2" "26.756.75"
Set-CellText "E2" "  -2.37%  "
Set-CellText "D3" "1.564.78"
Set-CellText "E3" "  +0.05%  "
Set-CellText "D5" "206.30"
Set-CellText "E5" "  -1.03%  "
Set-CellText "E6" "  -2.39%  "
Set-CellText "E7" "  -0.03%  "
Set-CellText "D8" "21.94"
Set-CellText "E8" "  -0.25%  "
Set-CellText "E9" "  -0.72%  "
Set-CellText "D10" "0.0584"
Set-CellText "E10" "  -1.33%  "
Set-CellText "E11" "  -0.63%  "
Set-CellText "D12" "1.786.65"
Set-CellText "E12" "  -0.14%  "
Set-CellText "D13" "1.566.76"
Set-CellText "E13" "  +0.02%  "
Set-CellText "E14" "  -2.51%  "
Set-CellText "E15" "  -0.54%  "
Set-CellText "D16" "26.807.45"
Set-CellText "D17" "61.48"
Set-CellText "D18" "213.85"
Set-CellText "E18" "  +0.59%  "
Set-CellText "E19" "  +1.27%  "
Set-CellText "D20" "0.0₃0676"
Set-CellText "E20" "  -1.87%  "
Set-CellText "E21" "  +0.10%  "
Set-CellText "E22" "  -0.49%  "
Set-CellText "D23" "9.32"
Set-CellText "E23" "  -1.96%  "
Set-CellText "D24" "2.01"
Set-CellText "E24" "  -0.45%  "
Set-CellText "D25" "153.17"
Set-CellText "E25" "  +0.26%  "
Set-CellText "E26" "  +0.68%  "
Set-CellText "E27" "  -0.35%  "
Set-CellText "E28" "  +0.08%  "
Set-CellText "E29" "  -1.35%  "
Set-CellText "E30" "  -1.35%  "
Set-CellText "E31" "  -3.73%  "
Set-CellText "E32" "  -1.75%  "
Set-CellText "D33" "1.384.03"
Set-CellText "E33" "  +0.77%  "
Set-CellText "E34" "  -1.28%  "
Set-CellText "E35" "  +0.59%  "
Set-CellText "E36" "  -0.86%  "
Set-CellText "D37" "0.923"
Set-CellText "E37" "  -4.09%  "
Set-CellText "E38" "  -2.61%  "
Set-CellText "E39" "  -1.62%  "
Set-CellText "D40" "0.814"
Set-CellText "E40" "  -0.75%  "
Set-CellText "E41" "  +0.07%  "
Set-CellText "D42" "0.990"
Set-CellText "E42" "  +1.60%  "
Set-CellText "D43" "5.35"
Set-CellText "E43" "  +1.77%  "
Set-CellText "B44" "MXToken"
Set-CellText "C44" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText "D44" "2.18"
Set-CellText "E44" "  +0.91%  "
Set-CellText "B45" "RenderToken"
Set-CellText "C45" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D45" "1.77"
Set-CellText "E45" "  -1.22%  "
Set-CellText "D46" "63.22"
Set-CellText "E46" "  -0.96%  "
Set-CellText "D47" "1.699.75"
Set-CellText "E47" "  -0.10%  "
Set-CellText "E48" "  -0.08%  "
Set-CellText "D49" "0.0₇0985"
Set-CellText "E49" "  -0.14%  "
Set-CellText "E50" "  -0.90%  "
Set-CellText "D51" "0.0492"
Set-CellText "E51" "  -0.64%  "
